# v0.2.1 (2019/9/8 23:25) Added Template.xlsx 方便进行统计
#
# The author filled in the previously-blank "elapsed seconds" helper
# formulas in column W (rows 3-10) and column AC (rows 4-10) of Sheet1,
# mirroring the pattern already used by the neighbouring K/Q columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (stand-alone, not part of the W4:W10 fill-down block)
$ws.Range("W3").Formula = "=IF(ISBLANK(T3),0,(MINUTE(V3)*60+SECOND(V3))-(MINUTE(U3)*60+SECOND(U3)))"

# Column W: rows 4-10 (fill down from the existing W5:W10 pattern)
for ($r = 4; $r -le 10; $r++) {
    $ws.Range("W$r").Formula = "=IF(ISBLANK(T$r),0,(MINUTE(V$r)*60+SECOND(V$r))-(MINUTE(U$r)*60+SECOND(U$r)))"
}

# Column AC: rows 4-10 (fill down from the existing AC5:AC7 pattern)
for ($r = 4; $r -le 10; $r++) {
    $ws.Range("AC$r").Formula = "=IF(ISBLANK(Z$r),0,(MINUTE(AB$r)*60+SECOND(AB$r))-(MINUTE(AA$r)*60+SECOND(AA$r)))"
}

# Reposition the view roughly where the author left it (best effort -
# the underlying window/pane model only tracks a single active selection).
$ws.Activate()
$ws.Range("F13").Select() | Out-Null

$wb.Save() | Out-Null
